# Trade #118 closed at 2026-02-17 09:28:25 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade (#118).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.68   # Current Capital
$summary.Range("B4").Value = 0.69      # Total P&L $
$summary.Range("B5").Value = 0.12      # Total P&L %
$summary.Range("B6").Value = 118       # Total Trades
$summary.Range("B7").Value = 53        # Winning Trades
$summary.Range("B9").Value = 44.92     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.68     # Capital
$status.Range("D4").Value = 118        # Trades
$status.Range("E4").Value = 0.69       # P&L $
$status.Range("F4").Value = 0.68       # P&L %
$status.Range("G4").Value = 44.92      # Win Rate %

# ---------------------------------------------------------------------------
# Append the new closed trade (#118) as row 119 on both the "All Trades"
# and "MarketMaking" sheets (they mirror each other).
# ---------------------------------------------------------------------------
$tradeSheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(119, 1).Value = 118

    # Column B looks like a date (YYYY-MM-DD) which Excel's automatic type
    # detection would otherwise convert into a date serial number. Force it
    # to be treated as plain text, write it, then drop the now-unneeded
    # explicit number format so the cell keeps the default style.
    $ws.Cells.Item(119, 2).NumberFormat = "@"
    $ws.Cells.Item(119, 2).Value = "2026-02-17"
    $ws.Cells.Item(119, 2).ClearFormats()

    $ws.Cells.Item(119, 3).Value = "09:28:18"
    $ws.Cells.Item(119, 4).Value = "MarketMaking"
    $ws.Cells.Item(119, 5).Value = "UP"
    $ws.Cells.Item(119, 6).Value = 0.15
    $ws.Cells.Item(119, 7).Value = 0.19
    $ws.Cells.Item(119, 8).Value = "CLOSED"
    $ws.Cells.Item(119, 9).Value = 26.6667
    $ws.Cells.Item(119, 10).Value = 0.04
    $ws.Cells.Item(119, 11).Value = 100.68
    $ws.Cells.Item(119, 12).Value = 0
    $ws.Cells.Item(119, 13).Value = 0
    $ws.Cells.Item(119, 14).Value = 0.6
    $ws.Cells.Item(119, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(119, 16).Value = "early_exit"
    $ws.Cells.Item(119, 17).Value = 0.15
}
